$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "18-25" placeholder text in column C with the actual
# target-market-budget formula (74,500,000 * 17%), written the way Excel
# itself would after typing the formula into C2 and then filling it down
# in successive chunks (which is why the saved file breaks the shared
# formula group into pieces of 64 rows).
$ws.Range("C2").Formula = "=74500000*0.17"
$ws.Range("C3:C66").Formula = "=74500000*0.17"
$ws.Range("C67:C130").Formula = "=74500000*0.17"
$ws.Range("C131:C194").Formula = "=74500000*0.17"
$ws.Range("C195:C230").Formula = "=74500000*0.17"

# Move the active selection to C2, matching where the author ended up
# after making the edit.
[void]$ws.Range("C2").Select()
